# Auto-generated Excel COM-interop script
# Commit: Add data for 2024-11-16
# Applies 176 individual cell value updates (YTD day-count increments)
# across "Citywide Totals", "By Neighborhood", and 29 neighborhood sheets.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 83
$ws.Range("H2").Value = 101
$ws.Range("J2").Value = 112
$ws.Range("F3").Value = 126
$ws.Range("G3").Value = 134
$ws.Range("J3").Value = 213
$ws.Range("K3").Value = 207
$ws.Range("B9").Value = 356
$ws.Range("E9").Value = 438
$ws.Range("H9").Value = 427
$ws.Range("J9").Value = 390
$ws.Range("K9").Value = 473
$ws.Range("B10").Value = 1256
$ws.Range("C10").Value = 1482
$ws.Range("D10").Value = 1696
$ws.Range("E10").Value = 2015
$ws.Range("F10").Value = 1990
$ws.Range("G10").Value = 865
$ws.Range("H10").Value = 565
$ws.Range("I10").Value = 807
$ws.Range("J10").Value = 684
$ws.Range("K10").Value = 631
$ws.Range("B11").Value = 1738
$ws.Range("C11").Value = 2089
$ws.Range("D11").Value = 2307
$ws.Range("E11").Value = 2672
$ws.Range("F11").Value = 2695
$ws.Range("G11").Value = 1509
$ws.Range("H11").Value = 1256
$ws.Range("I11").Value = 1610
$ws.Range("J11").Value = 1429
$ws.Range("K11").Value = 1474

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K5").Value = 22
$ws.Range("B8").Value = 66
$ws.Range("F8").Value = 135
$ws.Range("K8").Value = 63
$ws.Range("C9").Value = 6
$ws.Range("C10").Value = 11
$ws.Range("K10").Value = 15
$ws.Range("F21").Value = 24
$ws.Range("F27").Value = 26
$ws.Range("F28").Value = 116
$ws.Range("G28").Value = 82
$ws.Range("H28").Value = 71
$ws.Range("K28").Value = 86
$ws.Range("E29").Value = 24
$ws.Range("J29").Value = 21
$ws.Range("E32").Value = 141
$ws.Range("H32").Value = 73
$ws.Range("J36").Value = 58
$ws.Range("C41").Value = 27
$ws.Range("J41").Value = 32
$ws.Range("H48").Value = 7
$ws.Range("B53").Value = 246
$ws.Range("D53").Value = 560
$ws.Range("E53").Value = 671
$ws.Range("F53").Value = 590
$ws.Range("G53").Value = 239
$ws.Range("H53").Value = 190
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 2
$ws.Range("K62").Value = 21
$ws.Range("E65").Value = 45
$ws.Range("B68").Value = 12
$ws.Range("I70").Value = 32
$ws.Range("C74").Value = 39
$ws.Range("E76").Value = 92
$ws.Range("J76").Value = 42
$ws.Range("B77").Value = 77
$ws.Range("D78").Value = 65
$ws.Range("D83").Value = 28
$ws.Range("H83").Value = 21
$ws.Range("B89").Value = 23
$ws.Range("C90").Value = 4
$ws.Range("C92").Value = 24
$ws.Range("F95").Value = 58
$ws.Range("B99").Value = 1738
$ws.Range("C99").Value = 2089
$ws.Range("D99").Value = 2307
$ws.Range("E99").Value = 2672
$ws.Range("F99").Value = 2695
$ws.Range("G99").Value = 1509
$ws.Range("H99").Value = 1256
$ws.Range("I99").Value = 1610
$ws.Range("J99").Value = 1429
$ws.Range("K99").Value = 1474

# Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("E6").Value = 9
$ws.Range("J7").Value = 28
$ws.Range("E8").Value = 92
$ws.Range("J8").Value = 42

# Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("F6").Value = 17
$ws.Range("F7").Value = 26

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("B9").Value = 52
$ws.Range("B10").Value = 77

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K7").Value = 23
$ws.Range("B8").Value = 38
$ws.Range("F8").Value = 92
$ws.Range("B9").Value = 66
$ws.Range("F9").Value = 135
$ws.Range("K9").Value = 63

# Chinatown
$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("F3").Value = 2
$ws.Range("F9").Value = 24

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("H7").Value = 36
$ws.Range("E8").Value = 80
$ws.Range("E9").Value = 141
$ws.Range("H9").Value = 73

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 7
$ws.Range("J9").Value = 16
$ws.Range("J10").Value = 58

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("G3").Value = 20
$ws.Range("B9").Value = 199
$ws.Range("D9").Value = 492
$ws.Range("E9").Value = 591
$ws.Range("F9").Value = 518
$ws.Range("H9").Value = 98
$ws.Range("B10").Value = 246
$ws.Range("D10").Value = 560
$ws.Range("E10").Value = 671
$ws.Range("F10").Value = 590
$ws.Range("G10").Value = 239
$ws.Range("H10").Value = 190

# Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 22

# Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I7").Value = 14
$ws.Range("I8").Value = 32

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("E7").Value = 34
$ws.Range("E8").Value = 45

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("B4").Value = 6
$ws.Range("B6").Value = 23

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 7
$ws.Range("C6").Value = 18
$ws.Range("C7").Value = 27
$ws.Range("J7").Value = 32

# Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("D5").Value = 60
$ws.Range("D6").Value = 65

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("F2").Value = 6
$ws.Range("K3").Value = 17
$ws.Range("H7").Value = 25
$ws.Range("G8").Value = 37
$ws.Range("H8").Value = 27
$ws.Range("F9").Value = 116
$ws.Range("G9").Value = 82
$ws.Range("H9").Value = 71
$ws.Range("K9").Value = 86

# Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J7").Value = 5
$ws.Range("E8").Value = 17
$ws.Range("E9").Value = 24
$ws.Range("J9").Value = 21

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("H4").Value = 3
$ws.Range("D5").Value = 19
$ws.Range("D6").Value = 28
$ws.Range("H6").Value = 21

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("C6").Value = 33
$ws.Range("C7").Value = 39

# West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("C8").Value = 21
$ws.Range("C9").Value = 24

# Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K7").Value = 13
$ws.Range("K8").Value = 21

# Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 6

# Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 2
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 2

# West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("F6").Value = 51
$ws.Range("F7").Value = 58

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("H2").Value = 1
$ws.Range("H7").Value = 7

# O'Hare
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 12

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K5").Value = 2
$ws.Range("C6").Value = 9
$ws.Range("C7").Value = 11
$ws.Range("K7").Value = 15

# West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 4
